# PFAS_Rdata_Tissue.xlsx — "Add files via upload" re-edit
#
# The ID column (A2:A13) had three different bordered cell styles
# (a top/bottom-medium-ruled, left/right-thin box around the data
# block); the E and G columns ("KW_BW"/"HW_BW") were stored as raw
# fractions (e.g. 0.0114) instead of the per-mille values actually
# wanted (11.4). This pass strips the now-unwanted borders back to
# the plain "Normal" style (which also drops the artificial
# thick-bottom-border row height on every row), rescales the KW_BW /
# HW_BW columns by x1000, and leaves the selection where the author
# left off (M10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rescale columns E (KW_BW) and G (HW_BW), rows 2-13, by 1000 ---
$ws.Range("E2").Value = 11.4
$ws.Range("G2").Value = 4.7
$ws.Range("E3").Value = 10.8
$ws.Range("G3").Value = 4.0999999999999996
$ws.Range("E4").Value = 10.6
$ws.Range("G4").Value = 4.0999999999999996
$ws.Range("E5").Value = 10.6
$ws.Range("G5").Value = 4.4000000000000004
$ws.Range("E6").Value = 12
$ws.Range("G6").Value = 4.5
$ws.Range("E7").Value = 12.2
$ws.Range("G7").Value = 4.4000000000000004
$ws.Range("E8").Value = 13
$ws.Range("G8").Value = 5.4
$ws.Range("E9").Value = 11.3
$ws.Range("G9").Value = 4.2
$ws.Range("E10").Value = 13.9
$ws.Range("G10").Value = 4.3
$ws.Range("E11").Value = 13.2
$ws.Range("G11").Value = 4.2
$ws.Range("E12").Value = 13.9
$ws.Range("G12").Value = 4.2
$ws.Range("E13").Value = 14
$ws.Range("G13").Value = 4.3

# --- Strip the medium/thin box-border styling used on column A (ID) ---
# (reverts every cell in the used range back to the plain "Normal"
# style, removing the s="1"/"2"/"3" borders)
$ws.Range("A1:G13").Style = "Normal"

# --- Drop the explicit thick-bottom-border row height left behind by
#     the old borders, restoring the default row height on every row ---
$ws.Rows("1:13").AutoFit()

# --- Leave the selection where the author ended up ---
$ws.Range("M10").Select() | Out-Null
